$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The reaction is now between only two species (columns A and B);
# drop the old C1:Q1 values entirely so the sheet's used range shrinks.
$ws.Range("C1:Q1").ClearContents()

# Random, non-overlapping reaction indices for the remaining columns.
$ws.Range("A1").Value = 2
$ws.Range("B1").Value = 3
